$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Row 40: Ferienspass 8.8.24 - Popcorn
$row40 = $lo.ListRows.Add()
$ws.Range("A39").Copy()
$ws.Range("A40").PasteSpecial(-4122)
$ws.Range("A40").Value = 45512
$ws.Range("B40").Value = "Spez 1"
$ws.Range("C40").Value = "Popcorn"
$ws.Range("D39").Copy()
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("D40").Value = 5

# Row 41: Ferienspass 8.8.24 - Getraenke
$row41 = $lo.ListRows.Add()
$ws.Range("A39").Copy()
$ws.Range("A41").PasteSpecial(-4122)
$ws.Range("A41").Value = 45512
$ws.Range("B41").Value = "Spez 2"
$ws.Range("C41").Value = "Getränke"
$ws.Range("D39").Copy()
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("D41").Value = 3

$excel.CutCopyMode = 0

$ws.Range("C42").Select()
